# MITHEx_inputs.xlsx edit script
# Converts the "Plant Description" / "Optional Parameters" / "Input options"
# workbook from a single-primary/secondary-fluid tube/shell HX input sheet
# into a plate heat-exchanger (e-NTU / LMTD) input sheet.

$wb = $excel.ActiveWorkbook

$wsPlant = $wb.Worksheets.Item("Plant Description")
$wsOpt   = $wb.Worksheets.Item("Optional Parameters")
$wsOpts  = $wb.Worksheets.Item("Input options")

# ---------------------------------------------------------------------
# 1. "Input options" sheet - clear the old "Tube Pitch Type" column (C)
# ---------------------------------------------------------------------
# Select C9 here first (this is where the user last clicked on this sheet
# before moving back to "Plant Description").
$wsOpts.Range("C9").Select()

# Reformat C1:C3 to the plain/no-border style (copy from C4, which already
# uses that style) and then clear their text content.
$wsOpts.Range("C4").Copy()
$wsOpts.Range("C1:C3").PasteSpecial(-4122)
$wsOpts.Application.CutCopyMode = $false
$wsOpts.Range("C1:C3").ClearContents()

# ---------------------------------------------------------------------
# 2. "Plant Description" sheet - rework into primary/secondary PHE inputs
# ---------------------------------------------------------------------
$wsPlant.Columns.Item(1).ColumnWidth = 28.498697916666668

# Row 1: Thermal Power (kW) -> Thermal Power (MW), value 5
$wsPlant.Range("A1").Value2 = "Thermal Power (MW)"
$wsPlant.Range("B1").Value2 = 5

# Row 2: Primary Fluid / Sodium -- unchanged

# Row 3: Hot Temperature (C) -> Primary Hot Temperature (C), value 550
$wsPlant.Range("A3").Value2 = "Primary Hot Temperature (C)"
$wsPlant.Range("B3").Value2 = 550

# Row 4: Cold Temperature (C) -> Primary Cold Temperature (C), value 350
$wsPlant.Range("A4").Value2 = "Primary Cold Temperature (C)"
$wsPlant.Range("B4").Value2 = 350

# Row 5: Mass Flow Rate (kg/s) -> Primary Mass Flow Rate (kg/s)
$wsPlant.Range("A5").Value2 = "Primary Mass Flow Rate (kg/s)"

# Remove the old dropdown validation for Secondary Fluid that currently
# lives on B6 -- it is being replaced by a numeric "Primary Pressure" row,
# and the fluid-picker moves down to B7 (added back below).
$wsPlant.Range("B6").Validation.Delete()

# Row 6: Secondary Fluid -> Primary Pressure (kPa), value 100
# (give A6 the "highlighted / no top-bottom border" look used elsewhere
# for a sub-heading style row)
$wsPlant.Range("A6").Value2 = "Primary Pressure (kPa)"
$wsPlant.Range("B6").Value2 = 100
$wsPlant.Range("A6").Borders.Item(8).LineStyle = -4142
$wsPlant.Range("A6").Borders.Item(9).LineStyle = -4142

# Row 7: Hot Temperature (C) -> Secondary Fluid, value Air (with the
# dropdown list validation that used to be on B6)
$wsPlant.Range("A7").Value2 = "Secondary Fluid"
$wsPlant.Range("B7").Value2 = "Air"
$wsPlant.Range("B7").Validation.Add(3, 1, 1, "='Input options'!`$B`$2:`$B`$4")

# Row 8: Cold Temperature (C) -> Secondary Hot Temperature (C)
$wsPlant.Range("A8").Value2 = "Secondary Hot Temperature (C)"

# Row 9: Mass Flow Rate (kg/s) -> Secondary Cold Temperature (C), value 123
$wsPlant.Range("A9").Value2 = "Secondary Cold Temperature (C)"
$wsPlant.Range("B9").Value2 = 123

# Row 10 (new): Secondary Mass Flow Rate (kg/s), value 20
$wsPlant.Range("A9:B9").Copy()
$wsPlant.Range("A10:B10").PasteSpecial(-4122)
$wsPlant.Range("A11:B11").PasteSpecial(-4122)
$wsPlant.Application.CutCopyMode = $false

$wsPlant.Range("A10").Value2 = "Secondary Mass Flow Rate (kg/s)"
$wsPlant.Range("B10").Value2 = 20

# Row 11 (new): Secondary Pressure (kPa), value 2000
$wsPlant.Range("A11").Value2 = "Secondary Pressure (kPa)"
$wsPlant.Range("B11").Value2 = 2000

# ---------------------------------------------------------------------
# 3. "Optional Parameters" sheet - plate HX geometry inputs
# ---------------------------------------------------------------------
$wsOpt.Columns.Item(1).ColumnWidth = 25.998697916666668
$wsOpt.Columns.Item(2).ColumnWidth = 11.830729166666666

# Give A1/A2 the same "label" look used on Plant Description (green box)
$wsPlant.Range("A1").Copy()
$wsOpt.Range("A1:A2").PasteSpecial(-4122)
$wsPlant.Application.CutCopyMode = $false

# Row 1: Tube O.D. (in) -> Plate thickness (m), value 0.001
$wsOpt.Range("A1").Value2 = "Plate thickness (m)"
$wsOpt.Range("B1").Value2 = 0.001
$wsPlant.Range("B1").Copy()
$wsOpt.Range("B1").PasteSpecial(-4122)
$wsPlant.Application.CutCopyMode = $false
$wsOpt.Range("B1").Value2 = 0.001

# Row 2: Pitch type -> Plate material, value SS316
$wsOpt.Range("A2").Value2 = "Plate material"
$wsOpt.Range("B2").Value2 = "SS316"
$wsPlant.Range("B2").Copy()
$wsOpt.Range("B2").PasteSpecial(-4122)
$wsPlant.Application.CutCopyMode = $false
$wsOpt.Range("B2").Value2 = "SS316"

# Row 3 (new): Hydraulic Diameter (m), value 0.005
$wsOpts.Range("C3").Copy()
$wsOpt.Range("A3").PasteSpecial(-4122)
$wsOpts.Application.CutCopyMode = $false
$wsOpts.Range("B2").Copy()
$wsOpt.Range("B3").PasteSpecial(-4122)
$wsOpts.Application.CutCopyMode = $false
$wsOpt.Range("A3").Value2 = "Hydraulic Diameter (m)"
$wsOpt.Range("B3").Value2 = 0.005

# Row 4 (new): Primary Flow Velocity (m/s), value 0.5
$wsPlant.Range("A1:B1").Copy()
$wsOpt.Range("A4:B4").PasteSpecial(-4122)
$wsPlant.Application.CutCopyMode = $false
$wsOpt.Range("A4").Value2 = "Primary Flow Velocity (m/s)"
$wsOpt.Range("B4").Value2 = 0.5

# Row 5 (new): Secondary Flow Velocity (m/s), value 0.5
$wsPlant.Range("A1:B1").Copy()
$wsOpt.Range("A5:B5").PasteSpecial(-4122)
$wsPlant.Application.CutCopyMode = $false
$wsOpt.Range("A5").Value2 = "Secondary Flow Velocity (m/s)"
$wsOpt.Range("B5").Value2 = 0.5

# Rows 6-7 (new): blank spacer rows, formatted plain (same style as the
# cleared "Input options" C column)
$wsOpts.Range("C4").Copy()
$wsOpt.Range("A6:B7").PasteSpecial(-4122)
$wsOpts.Application.CutCopyMode = $false

$wsOpt.Range("B3").Select()

# ---------------------------------------------------------------------
# 4. Final selection / active sheet state
# ---------------------------------------------------------------------
$wsPlant.Activate()
$wsPlant.Range("B2").Select()
